$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 values (E2, F2, G2)
$ws.Range("E2").Value = 9500
$ws.Range("F2").Value = 20250
$ws.Range("G2").Value = -10750

# Add new row 4: Sausage, 150, kg, 55 (E4/F4/G4 left blank, like E3/F3/G3)
$ws.Range("A4").Value = "Sausage"
$ws.Range("B4").Value = 150
$ws.Range("C4").Value = "kg"
$ws.Range("D4").Value = 55

# E4/F4/G4 stay blank (mirroring E3/F3/G3); copy the blank cells from row 3
# so that the cells actually exist in the sheet rather than being entirely
# absent, matching the row 3 pattern of present-but-empty cells.
$ws.Range("E3:G3").Copy($ws.Range("E4:G4"))
